$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 627.75
$ws.Range("I107").Value = 585
$ws.Range("J107").Value = 642
$ws.Range("K107").Value = 585
$ws.Range("L107").Value = 642
$ws.Range("M107").Value = 1335
$ws.Range("N107").Value = -4482

$ws.Range("H129").Value = 2220.762
$ws.Range("I129").Value = 399.85715
$ws.Range("J129").Value = 3131.2144
$ws.Range("K129").Value = 1199.57145
$ws.Range("L129").Value = 9393.643199999999
$ws.Range("M129").Value = 3800.42855
$ws.Range("N129").Value = -19393.6432

$ws.Range("H132").Value = 2977786
$ws.Range("I132").Value = 3677647.2
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 11032941.6
$ws.Range("L132").Value = 10125
$ws.Range("M132").Value = -11030411.6
$ws.Range("N132").Value = -15185

$ws.Range("H137").Value = 946.6667
$ws.Range("I137").Value = 945.8461
$ws.Range("J137").Value = 952
$ws.Range("K137").Value = 2837.5383
$ws.Range("L137").Value = 2856
$ws.Range("M137").Value = -287.5383000000002
$ws.Range("N137").Value = -7956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 975.1579
$ws.Range("I2").Value = 851.86664
$ws.Range("J2").Value = 1437.5
$ws.Range("K2").Value = 851.86664
$ws.Range("L2").Value = 1437.5
$ws.Range("M2").Value = -738.86664
$ws.Range("N2").Value = -1663.5

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H32").Value = 6192.594
$ws.Range("I32").Value = 2815.125
$ws.Range("J32").Value = 20741.691
$ws.Range("K32").Value = 2815.125
$ws.Range("L32").Value = 20741.691
$ws.Range("M32").Value = -2528.125
$ws.Range("N32").Value = -21315.691

$ws.Range("H45").Value = 1673.2858
$ws.Range("I45").Value = 1483.1666
$ws.Range("J45").Value = 2814
$ws.Range("K45").Value = 1483.1666
$ws.Range("L45").Value = 2814
$ws.Range("M45").Value = -1106.1666
$ws.Range("N45").Value = -3568

$ws.Range("H61").Value = 1264.94
$ws.Range("I61").Value = 1154.1915
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1154.1915
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -942.1914999999999
$ws.Range("N61").Value = -3424

$ws.Range("H74").Value = 952
$ws.Range("I74").Value = 742.4
$ws.Range("J74").Value = 1047.2727
$ws.Range("K74").Value = 742.4
$ws.Range("L74").Value = 1047.2727
$ws.Range("M74").Value = 131.6
$ws.Range("N74").Value = -2795.2727

$ws.Range("H77").Value = 952
$ws.Range("I77").Value = 742.4
$ws.Range("J77").Value = 1047.2727
$ws.Range("K77").Value = 3712
$ws.Range("L77").Value = 5236.363499999999
$ws.Range("M77").Value = 656
$ws.Range("N77").Value = -13972.3635

$ws.Range("H116").Value = 975.1579
$ws.Range("I116").Value = 851.86664
$ws.Range("J116").Value = 1437.5
$ws.Range("K116").Value = 851.86664
$ws.Range("L116").Value = 1437.5
$ws.Range("M116").Value = 1442.13336
$ws.Range("N116").Value = -6025.5

$ws.Range("H122").Value = 1913.742
$ws.Range("I122").Value = 1325.6
$ws.Range("J122").Value = 2983.0908
$ws.Range("K122").Value = 3976.8
$ws.Range("L122").Value = 8949.2724
$ws.Range("M122").Value = -1526.8
$ws.Range("N122").Value = -13849.2724

$ws.Range("H132").Value = 1629.0212
$ws.Range("I132").Value = 1300.75
$ws.Range("J132").Value = 2703.3635
$ws.Range("K132").Value = 3902.25
$ws.Range("L132").Value = 8110.0905
$ws.Range("M132").Value = -1372.25
$ws.Range("N132").Value = -13170.0905

$ws.Range("H136").Value = 1264.94
$ws.Range("I136").Value = 1154.1915
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3462.5745
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -912.5744999999997
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 975.1579
$ws.Range("I3").Value = 851.86664
$ws.Range("J3").Value = 1437.5
$ws.Range("K3").Value = 851.86664
$ws.Range("L3").Value = 1437.5
$ws.Range("M3").Value = -737.86664
$ws.Range("N3").Value = -1665.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2421.5483
$ws.Range("I31").Value = 2001.5385
$ws.Range("J31").Value = 4605.6
$ws.Range("K31").Value = 2001.5385
$ws.Range("L31").Value = 4605.6
$ws.Range("M31").Value = -1706.5385
$ws.Range("N31").Value = -5195.6

$ws.Range("H34").Value = 2421.5483
$ws.Range("I34").Value = 2001.5385
$ws.Range("J34").Value = 4605.6
$ws.Range("K34").Value = 2001.5385
$ws.Range("L34").Value = 4605.6
$ws.Range("M34").Value = -1799.5385
$ws.Range("N34").Value = -5009.6

$ws.Range("H132").Value = 1766.1333
$ws.Range("I132").Value = 1513.6666
$ws.Range("J132").Value = 2355.2222
$ws.Range("K132").Value = 4540.9998
$ws.Range("L132").Value = 7065.6666
$ws.Range("M132").Value = -2010.9998
$ws.Range("N132").Value = -12125.6666

$ws.Range("H134").Value = 1260.0605
$ws.Range("I134").Value = 1078.76
$ws.Range("J134").Value = 1826.625
$ws.Range("K134").Value = 3236.28
$ws.Range("L134").Value = 5479.875
$ws.Range("M134").Value = -701.2799999999997
$ws.Range("N134").Value = -10549.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 732.44446
$ws.Range("I131").Value = 318.06668
$ws.Range("J131").Value = 1028.4286
$ws.Range("K131").Value = 954.2000400000001
$ws.Range("L131").Value = 3085.2858
$ws.Range("M131").Value = 4085.79996
$ws.Range("N131").Value = -13165.2858

$ws.Range("H132").Value = 1858.9286
$ws.Range("I132").Value = 1383
$ws.Range("J132").Value = 2123.3333
$ws.Range("K132").Value = 12447
$ws.Range("L132").Value = 19109.9997
$ws.Range("M132").Value = -9917
$ws.Range("N132").Value = -24169.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 44999.75
$ws.Range("J106").Value = 44999.75
$ws.Range("L106").Value = 44999.75
$ws.Range("N106").Value = -47523.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 372.90475
$ws.Range("I22").Value = 303.76923
$ws.Range("K22").Value = 303.76923
$ws.Range("M22").Value = -8.769229999999993

$ws.Range("H27").Value = 372.90475
$ws.Range("I27").Value = 303.76923
$ws.Range("K27").Value = 303.76923
$ws.Range("M27").Value = -196.76923

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H136").Value = 1973.2593
$ws.Range("I136").Value = 1144.4546
$ws.Range("J136").Value = 2543.0625
$ws.Range("K136").Value = 3433.3638
$ws.Range("L136").Value = 7629.1875
$ws.Range("M136").Value = -883.3638000000001
$ws.Range("N136").Value = -12729.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3978.9092
$ws.Range("J96").Value = 3989.8
$ws.Range("L96").Value = 3989.8
$ws.Range("N96").Value = -6735.8

$ws.Range("H97").Value = 22224.75
$ws.Range("J97").Value = 22224.75
$ws.Range("L97").Value = 22224.75
$ws.Range("N97").Value = -24206.75
